$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The paragraph that used to read "Versione1.27 del 05/05/2019" is
#    repurposed to announce the new version 1.28.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Versione1.27 del 05/05/2019", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Versione1.28 del 22/05/2019", 2)

# Locate that paragraph again (now holding the 1.28 text) so we can
# insert the new changelog entries right after it.
$rng = $d.Content
$rng.Find.Execute("Versione1.28 del 22/05/2019", $true, $true, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$verPara = $rng.Paragraphs(1)

# ------------------------------------------------------------------
# 2. Insert the two new bullet points describing the 1.28 release.
# ------------------------------------------------------------------
$anchor = $verPara.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()

$item1 = $d.Paragraphs($verPara.Range.Information(3) + 1).Range
$item1 = $rng.Paragraphs(1).Next().Range
$item1.Style = "Paragrafoelenco"
$item1.Text = "Aggiunta gestione Tipologie Interventi nel men$([char]0xF9) ""Tabelle/Tipologie Interventi"""
$item1.ListFormat.ApplyBulletDefault()

$item1After = $item1.Paragraphs(1).Range
$item1After.Collapse(0)
$item1After.InsertParagraphAfter()
$item2 = $item1.Paragraphs(1).Next().Range
$item2.Style = "Paragrafoelenco"
$item2.Text = "Aggiornata Lavorazione con gestione Dettaglio Economico e Piano Esterno (da completare)"
$item2.ListFormat.ApplyBulletDefault()

# ------------------------------------------------------------------
# 3. Blank paragraph, then a fresh "Versione1.27" sub-title paragraph
#    (recreating what used to be here before the new entries).
# ------------------------------------------------------------------
$item2After = $item2.Paragraphs(1).Range
$item2After.Collapse(0)
$item2After.InsertParagraphAfter()
$blank = $item2.Paragraphs(1).Next().Range
$blank.Style = "Normale"

$blankAfter = $blank.Paragraphs(1).Range
$blankAfter.Collapse(0)
$blankAfter.InsertParagraphAfter()
$verBack = $blank.Paragraphs(1).Next().Range
$verBack.Style = "Sottotitolo"
$verBack.Text = "Versione1.27 del 05/05/2019"

Write-Output "done"
